$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns B/C/D/E are treated as plain text so values like
# "30.647.63" or "0.4801" are not reinterpreted as numbers/dates by Excel.
$ws.Range('B2:E51').NumberFormat = '@'

$ws.Range('D2').Value = '30.647.63'
$ws.Range('E2').Value = '  +0.53%  '

$ws.Range('D3').Value = '1.949.45'
$ws.Range('E3').Value = '  +1.85%  '

$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').Value = '247.22'
$ws.Range('E5').Value = '  +1.06%  '

$ws.Range('E6').Value = '  +0.03%  '

$ws.Range('D7').Value = '0.4801'
$ws.Range('E7').Value = '  -0.63%  '

$ws.Range('E8').Value = '  +1.26%  '

$ws.Range('D9').Value = '0.06788'
$ws.Range('E9').Value = '  +1.00%  '

$ws.Range('D10').Value = '111.79'
$ws.Range('E10').Value = '  +1.89%  '

$ws.Range('E11').Value = '  +1.59%  '

$ws.Range('D12').Value = '1.960.96'
$ws.Range('E12').Value = '  +2.42%  '

$ws.Range('E13').Value = '  +1.70%  '

$ws.Range('D14').Value = '5.468'
$ws.Range('E14').Value = '  +3.81%  '

$ws.Range('D15').Value = '0.6835'
$ws.Range('E15').Value = '  +1.66%  '

$ws.Range('D16').Value = '293.22'
$ws.Range('E16').Value = '  +3.51%  '

$ws.Range('D17').Value = '30.662.71'
$ws.Range('E17').Value = '  +0.53%  '

$ws.Range('D18').Value = '13.18'
$ws.Range('E18').Value = '  +2.60%  '

$ws.Range('D19').Value = '5.645'
$ws.Range('E19').Value = '  +3.07%  '

$ws.Range('D20').Value = '0.000007661'
$ws.Range('E20').Value = '  +1.14%  '

$ws.Range('B21').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C21').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D21').Value = '2.205.16'
$ws.Range('E21').Value = '  +1.75%  '

$ws.Range('B22').Value = 'Dai'
$ws.Range('C22').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D22').Value = '0.9999'
$ws.Range('E22').Value = '  -0.01%  '

$ws.Range('D23').Value = '0.9997'
$ws.Range('E23').Value = '  +0.05%  '

$ws.Range('D24').Value = '6.583'
$ws.Range('E24').Value = '  +2.33%  '

$ws.Range('D25').Value = '9.732'

$ws.Range('D26').Value = '168.35'
$ws.Range('E26').Value = '  +2.37%  '

$ws.Range('D27').Value = '20.21'
$ws.Range('E27').Value = '  -0.08%  '

$ws.Range('D28').Value = '2.183'
$ws.Range('E28').Value = '  +2.94%  '

$ws.Range('D29').Value = '0.1087'
$ws.Range('E29').Value = '  +3.20%  '

$ws.Range('E30').Value = '  +2.04%  '

$ws.Range('D31').Value = '4.702'
$ws.Range('E31').Value = '  +16.49%  '

$ws.Range('D32').Value = '4.456'
$ws.Range('E32').Value = '  +6.98%  '

$ws.Range('D33').Value = '0.05055'
$ws.Range('E33').Value = '  +1.36%  '

$ws.Range('D34').Value = '0.7717'
$ws.Range('E34').Value = '  +5.49%  '

$ws.Range('D35').Value = '1.158'
$ws.Range('E35').Value = '  +2.17%  '

$ws.Range('D36').Value = '0.02071'
$ws.Range('E36').Value = '  +2.04%  '

$ws.Range('D37').Value = '2.735'
$ws.Range('E37').Value = '  +0.47%  '

$ws.Range('D38').Value = '2.697'
$ws.Range('E38').Value = '  +1.06%  '

$ws.Range('D39').Value = '2.053'
$ws.Range('E39').Value = '  +1.71%  '

$ws.Range('D40').Value = '110.58'
$ws.Range('E40').Value = '  -0.05%  '

$ws.Range('D41').Value = '0.4448'
$ws.Range('E41').Value = '  -0.19%  '

$ws.Range('D42').Value = '0.8702'
$ws.Range('E42').Value = '  +0.54%  '

$ws.Range('D43').Value = '5.967'
$ws.Range('E43').Value = '  +3.00%  '

$ws.Range('D44').Value = '0.9998'
$ws.Range('E44').Value = '  -0.02%  '

$ws.Range('D45').Value = '69.33'
$ws.Range('E45').Value = '  +1.89%  '

$ws.Range('D46').Value = '7.365'
$ws.Range('E46').Value = '  +0.46%  '

$ws.Range('D47').Value = '9.344'
$ws.Range('E47').Value = '  +1.01%  '

$ws.Range('D48').Value = '0.1250'
$ws.Range('E48').Value = '  +0.74%  '

$ws.Range('D49').Value = '47.91'
$ws.Range('E49').Value = '  -2.20%  '

$ws.Range('E50').Value = '  +2.30%  '

$ws.Range('D51').Value = '1.472'
$ws.Range('E51').Value = '  +1.45%  '
